$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "298.84"
Set-TextValue $ws.Range("E2") "-1.76%"
Set-TextValue $ws.Range("D3") "31.25"
Set-TextValue $ws.Range("E3") "-1.44%"
Set-TextValue $ws.Range("D4") "5.117"
Set-TextValue $ws.Range("E4") "-1.07%"
Set-TextValue $ws.Range("D5") "0.07946"
Set-TextValue $ws.Range("E5") "6.24%"
Set-TextValue $ws.Range("D6") "2.223"
Set-TextValue $ws.Range("E6") "-7.54%"
Set-TextValue $ws.Range("D7") "7.778"
Set-TextValue $ws.Range("E7") "-2.95%"
Set-TextValue $ws.Range("D8") "3.863"
Set-TextValue $ws.Range("E8") "-0.16%"
Set-TextValue $ws.Range("D9") "0.9263"
Set-TextValue $ws.Range("E9") "1.27%"
Set-TextValue $ws.Range("D10") "0.1730"
Set-TextValue $ws.Range("E10") "-0.36%"
Set-TextValue $ws.Range("D11") "0.07432"
Set-TextValue $ws.Range("E11") "-2.82%"
Set-TextValue $ws.Range("D12") "0.09438"
Set-TextValue $ws.Range("E12") "15.37%"
Set-TextValue $ws.Range("D13") "0.03032"
Set-TextValue $ws.Range("E13") "0.16%"
Set-TextValue $ws.Range("D14") "0.1004"
Set-TextValue $ws.Range("E14") "1.04%"
Set-TextValue $ws.Range("D15") "0.001515"
Set-TextValue $ws.Range("E15") "0.28%"
Set-TextValue $ws.Range("D16") "0.005858"
Set-TextValue $ws.Range("E16") "-4.74%"
Set-TextValue $ws.Range("D17") "3.480"
Set-TextValue $ws.Range("E17") "-0.47%"
Set-TextValue $ws.Range("E18") "1.41%"
Set-TextValue $ws.Range("E19") "0.28%"
Set-TextValue $ws.Range("D20") "0.1347"
Set-TextValue $ws.Range("E20") "1.36%"
Set-TextValue $ws.Range("D21") "3.921"
Set-TextValue $ws.Range("E21") "-15.80%"
Set-TextValue $ws.Range("D22") "0.1697"
Set-TextValue $ws.Range("E22") "8.43%"
Set-TextValue $ws.Range("D23") "0.04613"
Set-TextValue $ws.Range("E23") "-0.10%"
Set-TextValue $ws.Range("D24") "0.001244"
Set-TextValue $ws.Range("E24") "-1.50%"
Set-TextValue $ws.Range("D25") "0.004475"
Set-TextValue $ws.Range("E25") "-1.18%"
Set-TextValue $ws.Range("E26") "-7.78%"
Set-TextValue $ws.Range("D27") "0.0003392"
Set-TextValue $ws.Range("E27") "23.77%"
Set-TextValue $ws.Range("D39") "0.01765"
Set-TextValue $ws.Range("E39") "0.46%"
Set-TextValue $ws.Range("D40") "0.04605"
Set-TextValue $ws.Range("E40") "1.57%"
Set-TextValue $ws.Range("D41") "0.006975"
Set-TextValue $ws.Range("E41") "-5.47%"
Set-TextValue $ws.Range("D42") "0.1362"
Set-TextValue $ws.Range("E42") "-0.01%"
Set-TextValue $ws.Range("E43") "0.83%"
Set-TextValue $ws.Range("D44") "0.009566"
Set-TextValue $ws.Range("E44") "-12.12%"
Set-TextValue $ws.Range("D45") "0.00006283"
Set-TextValue $ws.Range("E45") "-1.69%"
Set-TextValue $ws.Range("D46") "0.00000000749"
Set-TextValue $ws.Range("E46") "-0.19%"
Set-TextValue $ws.Range("D47") "0.007966"
Set-TextValue $ws.Range("E47") "-19.42%"
Set-TextValue $ws.Range("D48") "0.7466"
Set-TextValue $ws.Range("E48") "-9.02%"
Set-TextValue $ws.Range("E49") "-0.19%"
Set-TextValue $ws.Range("E50") "-0.12%"
